$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update window position/size (workbookView)
$win = $excel.ActiveWindow
$win.Left = 11560
$win.Top = 300
$win.Width = 25600
$win.Height = 15520

# New compare-method column (E) and value column (F)
$ws.Range("E2").Value = "VQSR"
$ws.Range("F2").Value = 0.3

$ws.Range("E3").Value = "PASS"

$ws.Range("E4").Value = "FAIL"
$ws.Range("F4").Value = 0.1

$ws.Range("E5").Value = "PASS"
$ws.Range("F5").Value = 0.004

$ws.Range("E6").Value = "PASS"
$ws.Range("F6").Value = 0.0003

$ws.Range("E7").Value = "PASS"
$ws.Range("F7").Value = 0.0043

# Update selection to reflect the active cell after the edit
$ws.Range("F9").Select()
